$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last refreshed" timestamp note (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 05:10"

# --- Mexico (row 18): refreshed case counts ---
$ws.Range("D18").Value = 61871
$ws.Range("E18").Value = 15862

# --- Jamaica (row 135): refreshed case counts ---
$ws.Range("B135").Value = 581
$ws.Range("C135").Value = 6
$ws.Range("D135").Value = 290
$ws.Range("E135").Value = 282

# --- New country "Butan" inserted right after "Macao" (which sits at row 183).
#     This pushes San Martin (Parte Francesa), Puerto Rico, Eritrea and
#     Botsuana down by one row (184->185->186->187->188), and the newly
#     freed-up row 184 is populated with Butan's fresh data. ---
$ws.Range("A188").Value = "Botsuana"
$ws.Range("B188").Value = 35
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 20
$ws.Range("E188").Value = 14
$ws.Range("H188").Value = 1

$ws.Range("A187").Value = "Eritrea"
$ws.Range("B187").Value = 39
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 39
$ws.Range("E187").Value = 0
$ws.Range("H187").Value = 0

$ws.Range("A186").Value = "Puerto Rico"
$ws.Range("B186").Value = 39
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 1
$ws.Range("E186").Value = 36
$ws.Range("H186").Value = 2

$ws.Range("A185").Value = "San Martin (Parte Francesa)"
$ws.Range("B185").Value = 41
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 33
$ws.Range("E185").Value = 5
$ws.Range("H185").Value = 3

$ws.Range("A184").Value = "Butan"
$ws.Range("B184").Value = 43
$ws.Range("C184").Value = 10
$ws.Range("D184").Value = 6
$ws.Range("E184").Value = 37
$ws.Range("H184").Value = 0

# --- "Santa Lucia" moves above "Belice" in the ranking (rows 200/201
#     swap places along with their data). ---
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2
